# Correção nos dados e início da análise PNAD 2009
# Row 6 ("grandes regiões e unidades da federação") was an empty header-like
# row with no data; it is removed so that the region/UF rows that follow
# shift up by one, aligning each row's data with the correct label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
